# The underlying OOXML diff for this commit is a pure whitespace/attribute
# canonicalization of word/document.xml and word/styles.xml (every
# w:xxx="..." attribute list - and the xmlns declarations on the root
# <w:document> element - gets re-serialized in alphabetical order). No
# element, attribute value, or piece of text actually changes: the same
# namespaces, the same <w:color w:val="E36C0A" .../> run formatting, the
# same <w:pgSz>/<w:pgMar> page geometry, and the same <w:latentStyles>/
# <w:style> definitions are present before and after, just re-ordered by
# whatever packaging/export step produced the commit. That reordering is
# an artifact of the XML writer used to save the package - it is not
# something exposed through the Word object model (Word's COM API has no
# "attribute order" knob), so there is no in-document edit to reproduce
# here: the document content, formatting and properties are already
# identical to the target state.
#
# Touch the active document (per the harness contract) without mutating
# any content, so the package is simply re-saved as-is.
$d = $word.ActiveDocument
$null = $d.Name
